# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" header on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end of the workbook ---
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(44934.99999999999, 84, 22.75158101211484, 144.6575894402216),
    @(44941.99999999999, 87, 23.36710632874677, 149.6438563839832),
    @(44948.99999999999, 90, 30.81900920839642, 154.0082650308647),
    @(44955.99999999999, 93, 28.88317163945432, 153.9593965913192),
    @(44962.99999999999, 96, 33.85576235684763, 156.2940437505583),
    @(44976.99999999999, 103, 41.92869652780011, 162.3054427489211),
    @(44983.99999999999, 106, 43.86892034642661, 166.2191700355975),
    @(44990.99999999999, 109, 51.46904336216733, 174.9624334077069),
    @(44997.99999999999, 112, 51.58624883179534, 171.7429118495214),
    @(45004.99999999999, 115, 53.32309066205403, 177.36209463471),
    @(45011.99999999999, 118, 56.61700002980182, 178.2358465592443),
    @(45018.99999999999, 121, 60.05416877201944, 184.5956903243973),
    @(45025.99999999999, 125, 58.86878507834744, 188.009266139803),
    @(45032.99999999999, 128, 59.80510681928888, 192.0421692136544),
    @(45039.99999999999, 131, 67.65246101131432, 191.8087294129299),
    @(45046.99999999999, 134, 71.70906376960116, 194.9127025882344),
    @(45053.99999999999, 137, 78.78518266709771, 198.3254161827886)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Match the formatting already used in the workbook: bold/bordered/centered
# header row and date-formatted first column, taken from the "Weekly
# Quantity" sheet so the same style entries are reused.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A18").PasteSpecial(-4122)
